$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New attendance data: Name in column A, time-of-day (with leading space) in
# column B, column C left blank. Row order reshuffled per the source data,
# and the former "JENSEN HUANG" row is dropped entirely, shrinking the used
# range from A1:C10 to A1:C9.
$names = @(
    "DENZEL WASHINGTON",
    "CHARLES LECLERC",
    "TOM CRUISE",
    "JEFF BESOS",
    "LEO MESSI",
    "NELSON MANDELA",
    "SAM ALTMAN",
    "STEVE JOBS"
)
$times = @(
    " 01:44:04",
    " 01:44:59",
    " 02:15:22",
    " 02:16:02",
    " 02:16:19",
    " 02:16:35",
    " 02:17:05",
    " 02:27:45"
)

for ($i = 0; $i -lt $names.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $names[$i]
    $ws.Cells.Item($row, 2).Value = $times[$i]
    $ws.Cells.Item($row, 3).Value = ""
}

# Drop the now-unused last row (old row 10) so the sheet's used range
# shrinks from A1:C10 down to A1:C9.
$ws.Rows.Item(10).Delete()
